# Weekly price-sheet update: insert a new Cebollín record for "La Araucanía"
# (Vega Modelo de Temuco) at row 646, pushing the existing rows 646-765 down
# to 647-766 and extending the used range to A1:R766.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 646 -- this shifts rows 646..765 down to 647..766
# and (matching native Excel behaviour) copies the formatting of the row
# above, so column D keeps its date number format.
$ws.Rows.Item(646).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(646, 1).Value  = 10
$ws.Cells.Item(646, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(646, 3).Value  = "La Araucanía"
$ws.Cells.Item(646, 4).Value  = 45275
$ws.Cells.Item(646, 5).Value  = 9
$ws.Cells.Item(646, 6).Value  = 100112037
$ws.Cells.Item(646, 7).Value  = "Cebollín"
$ws.Cells.Item(646, 8).Value  = "Sin especificar"
$ws.Cells.Item(646, 9).Value  = "Primera"
$ws.Cells.Item(646, 10).Value = 155
$ws.Cells.Item(646, 11).Value = 7000
$ws.Cells.Item(646, 12).Value = 7000
$ws.Cells.Item(646, 13).Value = 7000
$ws.Cells.Item(646, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(646, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(646, 16).Value = 583
$ws.Cells.Item(646, 17).Value = 12
$ws.Cells.Item(646, 18).Value = "Hortaliza"
